$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BaseRate_Fields")

# Add new row 6, mirroring the layout/format of row 5 (the prior last row)
$ws.Range("A5:E5").Copy()
$ws.Range("A6:E6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# rowid -> "5" stored as text (matches existing rowid column pattern)
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "5"
$ws.Range("B6").Value = "01_TL_Base_Rates_27122019"
$ws.Range("E6").Value = "FINASTRA_CCB_BASERATE_SY_GROUP1_20190921.csv"

# New hyperlinks on C6 / D6, same external target & display text as the other rows
$ws.Hyperlinks.Add($ws.Range("D6"), "\\DataSet\TL_DataSet\BaseRates_GSFile\", "", "", "\DataSet\TL_DataSet\BaseRates_GSFile\")
$ws.Hyperlinks.Add($ws.Range("C6"), "\\DataSet\TL_DataSet\BaseRates_GSFile\", "", "", "\DataSet\TL_DataSet\BaseRates_GSFile\")

# Restore text format + the real path values (Hyperlinks.Add overwrote cell text with the display text)
$ws.Range("C6:D6").NumberFormat = "@"
$ws.Range("C6").Value = "\DataSet\NewUATDeals_DataSet\Transformation_Layer\TL_Base_Rate\"
$ws.Range("D6").Value = "\DataSet\NewUATDeals_DataSet\Transformation_Layer\TL_Base_Rate\BaseRates_Files_27DEC2019\"

# Update the active selection
$ws.Range("C9").Select()
